$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a handful of "Regalo" descriptions (drop apostrophes / shorten names)
$ws.Range("B5").Value = "Pacchetto Office"
$ws.Range("B34").Value = "Il gioco dell oca "
$ws.Range("B42").Value = "Enigmista"
$ws.Range("B49").Value = "Allegro chirurgo"
$ws.Range("B84").Value = "Assassin Creed"
$ws.Range("B91").Value = "Marvel Avengers"
$ws.Range("B105").Value = "Luigi Mansion"

# Leave the selection on the first cell that was edited, scrolled back to the top
$ws.Range("B5").Select()
